$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/number-format on D-column cells whose new values are unambiguous
# numeric-looking strings, so Excel keeps them as text (matching the source data,
# which stores every Price cell as a string even when it looks like a plain number).
$textForceCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D14", "D18", "D20", "D22", "D23", "D24", "D26", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D43", "D44", "D46", "D49", "D50")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.376.19'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.363.37'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '310.54'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '104.43'
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('D7').Value = '0.527'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('D10').Value = '36.25'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').Value = '52.94'
$ws.Range('E11').Value = '  +1.77%  '
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').Value = '7.02'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '2.726.29'
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('E16').Value = '  +6.01%  '
$ws.Range('D17').Value = '2.349.92'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').Value = '0.814'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').Value = '43.377.93'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').Value = '12.00'
$ws.Range('E20').Value = '  -4.32%  '
$ws.Range('D21').Value = '0.0₃0930'
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('D22').Value = '6.28'
$ws.Range('E22').Value = '  +3.55%  '
$ws.Range('D23').Value = '68.38'
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').Value = '243.30'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('E25').Value = '  +2.46%  '
$ws.Range('D26').Value = '2.63'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = '26.05'
$ws.Range('E28').Value = '  +9.07%  '
$ws.Range('E29').Value = '  +8.82%  '
$ws.Range('D30').Value = '36.59'
$ws.Range('E30').Value = '  -5.84%  '
$ws.Range('D31').Value = '9.62'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').Value = '162.42'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').Value = '5.31'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').Value = '18.38'
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('E36').Value = '  +6.42%  '
$ws.Range('D37').Value = '3.12'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').Value = '0.0741'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.95'
$ws.Range('E39').Value = '  +6.44%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '4.63'
$ws.Range('E40').Value = '  +10.55%  '
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').Value = '2.42'
$ws.Range('E43').Value = '  +4.65%  '
$ws.Range('D44').Value = '19.83'
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.997.78'
$ws.Range('E45').Value = '  +2.18%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0292'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  +2.39%  '
$ws.Range('E48').Value = '  +6.01%  '
$ws.Range('D49').Value = '58.18'
$ws.Range('E49').Value = '  +5.79%  '
$ws.Range('D50').Value = '2.91'
$ws.Range('E50').Value = '  -3.48%  '
$ws.Range('E51').Value = '  +2.86%  '
